# Generate Report for Archive
# - Update the localization status from "Ready for handoff" to "In Translation"
#   on all three sheets (Overview, zh-cn, de-de).
# - Re-fit the "status" columns that held that text to their new (narrower)
#   content width, mirroring Excel's behaviour when a report is regenerated.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: status appears in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: status is in column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: status is in column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
